$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text label change ---
# G20 currently reads "cache line" -> append "(個數)"
$ws.Range("G20").Value = "cache line(個數)"

# --- Update underlying data values (K4, J9:K12, J14:J17) ---
$ws.Range("K4").Value = 18273

$ws.Range("J9").Value = 2048
$ws.Range("K9").Value = 1024

$ws.Range("J10").Value = 1536
$ws.Range("K10").Value = 768

$ws.Range("J11").Value = 3072
$ws.Range("K11").Value = 1536

$ws.Range("J12").Value = 2048
$ws.Range("K12").Value = 1024

$ws.Range("J14").Value = 65536
$ws.Range("J15").Value = 49152
$ws.Range("J16").Value = 98304
$ws.Range("J17").Value = 65536

# --- New CEILING formulas in column K for rows 20-23 ---
$ws.Range("K20").Formula = "=CEILING(J20,1)"
$ws.Range("K21:K23").Formula = "=CEILING(J21,1)"

# --- Update selection to match the authored state ---
$ws.Range("E7").Select()
